# Generate Report for Handoff
# Adds two new handed-off files (6be9e9f0-3afe-4b08-bed9-7db6a91b7300 and
# eff29242-06e4-4f2c-b221-51cf7ed86c65) as new rows 4 & 5 on the Overview,
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# -------------------------------------------------------------------------
# Overview sheet (sheet1) - rows 4 and 5
# -------------------------------------------------------------------------

$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6be9e9f03afe4b08bed97db6a91b7300/e2e/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md", "", "", "6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md") | Out-Null
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-39-13 14:39:56"

$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/eff2924206e44f2cb22151cf7ed86c65/e2e/eff29242-06e4-4f2c-b221-51cf7ed86c65.md", "", "", "eff29242-06e4-4f2c-b221-51cf7ed86c65.md") | Out-Null
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-39-13 14:39:56"

# -------------------------------------------------------------------------
# zh-cn sheet (sheet2) - rows 4 and 5
# -------------------------------------------------------------------------

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6be9e9f03afe4b08bed97db6a91b7300/e2e/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md", "", "", "6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/6be9e9f03afe4b08bed97db6a91b7300/e2e/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md", "", "", ".md") | Out-Null
$zh.Range("C4").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9671892dd72fe1969c37969b6fcf905d20ecf4d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.89587c9ce092ae372c4895803a38c93c78271ad8.zh-cn.xlf", "", "", "6be9e9f0-3afe-4b08-bed9-7db6a91b7300.89587c9ce092ae372c4895803a38c93c78271ad8.zh-cn.xlf") | Out-Null
$zh.Range("E4").Value = "2016-03-13 14:39:53"
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/eff2924206e44f2cb22151cf7ed86c65/e2e/eff29242-06e4-4f2c-b221-51cf7ed86c65.md", "", "", "eff29242-06e4-4f2c-b221-51cf7ed86c65.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/eff2924206e44f2cb22151cf7ed86c65/e2e/eff29242-06e4-4f2c-b221-51cf7ed86c65.md", "", "", ".md") | Out-Null
$zh.Range("C5").Value = "Ready for handoff"
$zh.Hyperlinks.Add($zh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3212d225e8d7a9f3328a85b98493441103935ea1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eff29242-06e4-4f2c-b221-51cf7ed86c65.1d65c7b145e76c7c4a49d2c973f41be13bd32174.zh-cn.xlf", "", "", "eff29242-06e4-4f2c-b221-51cf7ed86c65.1d65c7b145e76c7c4a49d2c973f41be13bd32174.zh-cn.xlf") | Out-Null
$zh.Range("E5").Value = "2016-03-13 14:39:53"
$zh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H5").Value = "0001-01-01 00:00:00"
$zh.Range("I5").Value = "Include"

# -------------------------------------------------------------------------
# de-de sheet (sheet3) - rows 4 and 5
# -------------------------------------------------------------------------

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6be9e9f03afe4b08bed97db6a91b7300/e2e/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md", "", "", "6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/6be9e9f03afe4b08bed97db6a91b7300/e2e/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.md", "", "", ".md") | Out-Null
$de.Range("C4").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a9a06a79d856b28feaef0683b7af86b305cf0f5d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6be9e9f0-3afe-4b08-bed9-7db6a91b7300.89587c9ce092ae372c4895803a38c93c78271ad8.de-de.xlf", "", "", "6be9e9f0-3afe-4b08-bed9-7db6a91b7300.89587c9ce092ae372c4895803a38c93c78271ad8.de-de.xlf") | Out-Null
$de.Range("E4").Value = "2016-03-13 14:39:56"
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "Include"

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/eff2924206e44f2cb22151cf7ed86c65/e2e/eff29242-06e4-4f2c-b221-51cf7ed86c65.md", "", "", "eff29242-06e4-4f2c-b221-51cf7ed86c65.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/eff2924206e44f2cb22151cf7ed86c65/e2e/eff29242-06e4-4f2c-b221-51cf7ed86c65.md", "", "", ".md") | Out-Null
$de.Range("C5").Value = "Ready for handoff"
$de.Hyperlinks.Add($de.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8576cc1a64bc3c4f12c96f4d5ebba03552af82ec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eff29242-06e4-4f2c-b221-51cf7ed86c65.1d65c7b145e76c7c4a49d2c973f41be13bd32174.de-de.xlf", "", "", "eff29242-06e4-4f2c-b221-51cf7ed86c65.1d65c7b145e76c7c4a49d2c973f41be13bd32174.de-de.xlf") | Out-Null
$de.Range("E5").Value = "2016-03-13 14:39:56"
$de.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H5").Value = "0001-01-01 00:00:00"
$de.Range("I5").Value = "Include"
